# RPA datasets push 2024-07-18
# Row 3 (엠83) gets a fresh demand-forecast entry; rows 4-11 shift their
# B/C/E/F (date / price band / amount / underwriter) values up from the
# row immediately below them (old row N+1 -> new row N); row 19's
# confirmed price (column D) is finalized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data that lands on row 3 (엠83)
$ws.Cells.Item(3, 2).Value = "2024.08.02~08.08"
$ws.Cells.Item(3, 3).Value = "13,600~15,300"
$ws.Cells.Item(3, 5).Value = 21308
$ws.Cells.Item(3, 6).Value = "NH투자증권"

# Rows 4-11: pull the B/C/E/F values up from the row below (as captured
# from the sheet before this edit)
$ws.Cells.Item(4, 2).Value = "2024.08.01~08.07"
$ws.Cells.Item(4, 3).Value = "11,000~13,000"
$ws.Cells.Item(4, 5).Value = 16500
$ws.Cells.Item(4, 6).Value = "신영증권,유진투자증권"

$ws.Cells.Item(5, 2).Value = "2024.07.31~08.06"
$ws.Cells.Item(5, 3).Value = "32,000~40,200"
$ws.Cells.Item(5, 5).Value = 78720
$ws.Cells.Item(5, 6).Value = "삼성증권"

$ws.Cells.Item(6, 2).Value = "2024.07.30~08.05"
$ws.Cells.Item(6, 3).Value = "12,500~15,500"
$ws.Cells.Item(6, 5).Value = 17500
$ws.Cells.Item(6, 6).Value = "하나증권"

$ws.Cells.Item(7, 2).Value = "2024.07.30~08.05"
$ws.Cells.Item(7, 3).Value = "13,800~15,700"
$ws.Cells.Item(7, 5).Value = 42471
$ws.Cells.Item(7, 6).Value = "미래에셋증권"

$ws.Cells.Item(8, 2).Value = "2024.07.29~07.30"
$ws.Cells.Item(8, 3).Value = "2,000~2,000"
$ws.Cells.Item(8, 5).Value = 11600
$ws.Cells.Item(8, 6).Value = "교보증권"

$ws.Cells.Item(9, 2).Value = "2024.07.29~08.02"
$ws.Cells.Item(9, 3).Value = "24,000~29,000"
$ws.Cells.Item(9, 5).Value = 24000
$ws.Cells.Item(9, 6).Value = "한국투자증권"

$ws.Cells.Item(10, 2).Value = "2024.07.29~08.02"
$ws.Cells.Item(10, 3).Value = "4,900~5,700"
$ws.Cells.Item(10, 5).Value = 14700
$ws.Cells.Item(10, 6).Value = "대신증권"

$ws.Cells.Item(11, 2).Value = "2024.07.29~08.02"
$ws.Cells.Item(11, 3).Value = "18,000~21,000"
$ws.Cells.Item(11, 5).Value = 13518
$ws.Cells.Item(11, 6).Value = "키움증권"

# Row 19 (산일전기(유가)): confirmed price is now set
$ws.Cells.Item(19, 4).Value = "35000"
